# Mississippi overview workbook - convert numeric "count" cells to literal
# text strings (matching upstream COMM export behaviour) and append a
# "Total" row to the County sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overall": A2 numeric 545 -> text "545"
# ---------------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("A2").Value = "'" + $wsOverall.Range("A2").Value2

# ---------------------------------------------------------------------------
# Sheet "County": B2:B71 numeric -> text (same values), rows 72-75 get
# re-expressed as percentages/currency text, and a new Total row (76) is
# appended.
# ---------------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

for ($r = 2; $r -le 71; $r++) {
    $cell = $wsCounty.Cells.Item($r, 2)
    $cell.Value = "'" + $cell.Value2
}

foreach ($r in 72, 73, 74, 75) {
    $wsCounty.Cells.Item($r, 2).Value = "'0.00%"
    $wsCounty.Cells.Item($r, 3).Value = "'`$0"
    $wsCounty.Cells.Item($r, 4).Value = "'0.00%"
    $wsCounty.Cells.Item($r, 5).Value = "'0.00%"
    $wsCounty.Cells.Item($r, 6).Value = "'0.00%"
}

$wsCounty.Range("A76").Value = "'Total"
$wsCounty.Range("B76").Value = "'545"
$wsCounty.Range("C76").Value = "'`$855,145,071"
$wsCounty.Range("D76").Value = "'6.46%"
$wsCounty.Range("E76").Value = "'-26.56%"
$wsCounty.Range("F76").Value = "'74.13%"

# ---------------------------------------------------------------------------
# Sheet "Congressional District": B2:B6 numeric -> text (same values)
# ---------------------------------------------------------------------------
$wsDistrict = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 6; $r++) {
    $cell = $wsDistrict.Cells.Item($r, 2)
    $cell.Value = "'" + $cell.Value2
}

# ---------------------------------------------------------------------------
# Sheet "Size": B2:B8 numeric -> text (same values)
# ---------------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 8; $r++) {
    $cell = $wsSize.Cells.Item($r, 2)
    $cell.Value = "'" + $cell.Value2
}

# ---------------------------------------------------------------------------
# Sheet "Subsector": B2:B14 numeric -> text (same values)
# ---------------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 14; $r++) {
    $cell = $wsSubsector.Cells.Item($r, 2)
    $cell.Value = "'" + $cell.Value2
}
